# Apply the crypto-price refresh captured in the commit diff.
# Cell values in columns D/E are stored as literal TEXT in the workbook
# (not numbers), so values that look numeric ("591.41", "0.0000164", ...)
# must be forced to text -- otherwise Excel auto-converts them to the
# Number type on assignment. We do this the same way a user would in the
# UI: type a leading apostrophe, then clear the resulting "Text" quote-
# prefix formatting so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '66.777.54'
$ws.Range("E2").Value = '  +1.17%  '
# Row 3
$ws.Range("D3").Value = '3.319.55'
$ws.Range("E3").Value = '  +0.50%  '
# Row 4
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").Value = "'591.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.95%  '
# Row 6
$ws.Range("D6").Value = "'182.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.35%  '
# Row 7
$ws.Range("D7").Value = "'0.640"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.37%  '
# Row 8
$ws.Range("E8").Value = '  +0.04%  '
# Row 9
$ws.Range("D9").Value = '3.318.64'
$ws.Range("E9").Value = '  +0.48%  '
# Row 10
$ws.Range("D10").Value = "'0.127"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.52%  '
# Row 11
$ws.Range("E11").Value = '  +2.74%  '
# Row 12
$ws.Range("D12").Value = "'0.405"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.80%  '
# Row 13
$ws.Range("D13").Value = '3.892.38'
$ws.Range("E13").Value = '  +0.58%  '
# Row 14
$ws.Range("E14").Value = '  -2.23%  '
# Row 15
$ws.Range("D15").Value = '66.773.63'
$ws.Range("E15").Value = '  +1.05%  '
# Row 16
$ws.Range("D16").Value = "'26.78"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.35%  '
# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.324.74'
$ws.Range("E17").Value = '  +0.55%  '
# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").Value = "'0.0000164"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.14%  '
# Row 19
$ws.Range("D19").Value = "'430.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.62%  '
# Row 20
$ws.Range("D20").Value = "'5.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.60%  '
# Row 21
$ws.Range("D21").Value = "'13.10"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.73%  '
# Row 22
$ws.Range("D22").Value = "'7.34"
$ws.Range("D22").ClearFormats()
# Row 23
$ws.Range("E23").Value = '  +0.26%  '
# Row 24
$ws.Range("D24").Value = "'71.62"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.41%  '
# Row 25
$ws.Range("D25").Value = "'5.75"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.31%  '
# Row 26
$ws.Range("D26").Value = '3.450.37'
$ws.Range("E26").Value = '  +0.52%  '
# Row 27
$ws.Range("D27").Value = "'0.515"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.29%  '
# Row 28
$ws.Range("E28").Value = '  +6.34%  '
# Row 29
$ws.Range("E29").Value = '  +0.61%  '
# Row 30
$ws.Range("D30").Value = "'9.29"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.50%  '
# Row 31
$ws.Range("E31").Value = '  -0.17%  '
# Row 32
$ws.Range("E32").Value = '  -0.66%  '
# Row 33
$ws.Range("D33").Value = "'22.43"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.14%  '
# Row 34
$ws.Range("E34").Value = '  +0.11%  '
# Row 35
$ws.Range("D35").Value = "'5.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.25%  '
# Row 36
$ws.Range("D36").Value = "'6.61"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.16%  '
# Row 37
$ws.Range("E37").Value = '  -0.54%  '
# Row 38
$ws.Range("D38").Value = "'159.13"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.04%  '
# Row 39
$ws.Range("D39").Value = "'1.44"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.58%  '
# Row 40
$ws.Range("D40").Value = '2.872.84'
$ws.Range("E40").Value = '  +3.27%  '
# Row 41
$ws.Range("D41").Value = "'1.80"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.32%  '
# Row 42
$ws.Range("D42").Value = "'26.49"
$ws.Range("D42").ClearFormats()
# Row 43
$ws.Range("D43").Value = "'4.36"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.21%  '
# Row 44
$ws.Range("D44").Value = "'0.754"
$ws.Range("D44").ClearFormats()
# Row 45
$ws.Range("D45").Value = "'39.79"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.41%  '
# Row 46
$ws.Range("D46").Value = "'5.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.08%  '
# Row 47
$ws.Range("E47").Value = '  +0.98%  '
# Row 48
$ws.Range("D48").Value = "'0.0644"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.67%  '
# Row 49
$ws.Range("D49").Value = "'313.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.69%  '
# Row 50
$ws.Range("D50").Value = "'23.06"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.22%  '
# Row 51
$ws.Range("D51").Value = "'0.0272"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.13%  '
